$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-23 Tuesday" "2024-01-24 Wednesday"

Replace-Text "816÷7=116, 4" "686÷2=343, 0"
Replace-Text "141÷9=15, 6" "278÷7=39, 5"
Replace-Text "840÷3=280, 0" "660÷4=165, 0"
Replace-Text "599÷4=149, 3" "279÷4=69, 3"
Replace-Text "421÷7=60, 1" "399÷9=44, 3"

Replace-Text "396÷7=56, 4" "401÷9=44, 5"
Replace-Text "224÷5=44, 4" "496÷7=70, 6"
Replace-Text "364÷6=60, 4" "694÷7=99, 1"
Replace-Text "530÷3=176, 2" "430÷7=61, 3"
Replace-Text "128÷8=16, 0" "524÷8=65, 4"

Replace-Text "182÷9=20, 2" "169÷2=84, 1"
Replace-Text "307÷7=43, 6" "737÷9=81, 8"
Replace-Text "353÷3=117, 2" "186÷4=46, 2"
Replace-Text "316÷5=63, 1" "430÷2=215, 0"
Replace-Text "478÷9=53, 1" "314÷9=34, 8"

Replace-Text "542÷8=67, 6" "847÷7=121, 0"
Replace-Text "253÷6=42, 1" "808÷4=202, 0"
Replace-Text "498÷6=83, 0" "318÷6=53, 0"
Replace-Text "647÷2=323, 1" "639÷2=319, 1"
Replace-Text "634÷6=105, 4" "976÷7=139, 3"

Replace-Text "489÷5=97, 4" "568÷5=113, 3"
Replace-Text "950÷2=475, 0" "622÷2=311, 0"
Replace-Text "302÷5=60, 2" "833÷8=104, 1"
Replace-Text "351÷3=117, 0" "630÷4=157, 2"
Replace-Text "901÷2=450, 1" "692÷5=138, 2"

Write-Output "Done applying replacements"
